# Sintaxe_Exemplo.xlsx - apply the commit's edits via Excel COM interop.
# Commit message: "Update para executar com o BD_CODIGOS e manter a ordenacao
# correta dos labels"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update the explanatory cell-comments text (authored by Rayner Santos)
#    Each comment keeps its "Rayner Santos:" bold header line followed by
#    the (revised) body text.
# ---------------------------------------------------------------------

$ws.Range("B1").Comment.Text('Rayner Santos:' + [char]10 + 'Informar as colunas (bandeiras) que representa as colunas da tabela, separados por ", " (virgula e um espaço).')

$ws.Range("C1").Comment.Text('Rayner Santos:' + [char]10 + 'Informar o nome do cabeçalho desejado para as colunas da bandeira. Coloque o cabeçalho separado por ", " (virgula e um espaço).')

$ws.Range("D1").Comment.Text('Rayner Santos:' + [char]10 + 'Nome da variável que representa a linha da tabela. ' + [char]10 + 'Obs.: quando a variável que representa a linha for de uma tabela MULTIPLA é necessário colocar um nome diferente dos Valores_Agrup.' + [char]10 + 'Exemplo: Variáveis REC_1, REC_2, REC_3 que referem a uma MULTIPLA, o nome da variável que representa a linha da tabela poderá ser REC, ou seja, diferente dos nomes das variáveis que serão preenchidas no campo: “Valores_Agrup”.')

$ws.Range("E1").Comment.Text('Rayner Santos:' + [char]10 + 'A tabela deverá contabilizar os casos de `NS/NR`? Escreva SIM ou NAO (em maiuscula e sem acento). ')

$ws.Range("F1").Comment.Text('Rayner Santos:' + [char]10 + 'Somente quando a variável que representará as linhas for do TipoTabela IPA_10 ou IPA_5, informar os valores a serem considerados para o BTB separados por ", " (virgula e um espaço).')

$ws.Range("G1").Comment.Text('Rayner Santos:' + [char]10 + 'Somente quando a variável que representará as linhas for do TipoTabela IPA_10 ou IPA_5, informar os valores a serem considerados para o TTB separados por ", " (virgula e um espaço).')

$ws.Range("H1").Comment.Text('Rayner Santos:' + [char]10 + 'Informar somente quando a variável na linha for MULTIPLA (deverá informar o nome das colunas que descrevem a variável múltipla. Coloque  separdo por ", " (virgula e um espaço).')

# ---------------------------------------------------------------------
# 2) Fill in the missing BTB/TTB example values on the IPA_5 example row
#    (row 6) and clear the old free-text example that used to sit in the
#    "Valores_Agrup" column for that row.
# ---------------------------------------------------------------------

$ws.Range("F6").Value2 = "1, 2"
$ws.Range("G6").Value2 = "4, 5"
$ws.Range("H6").Value2 = ""

# ---------------------------------------------------------------------
# 3) Resize a few columns so the sheet reads comfortably with the revised
#    comments/content (values taken from the saved workbook; the inputs
#    below are pre-compensated for the ColumnWidth -> stored-width
#    rounding/offset so the saved width lands on the intended value).
# ---------------------------------------------------------------------

$ws.Columns.Item(2).ColumnWidth = 35.833333333333336
$ws.Columns.Item(3).ColumnWidth = 54.166666666666664
$ws.Columns.Item(6).ColumnWidth = 5.666666666666667
$ws.Columns.Item(7).ColumnWidth = 4.666666666666667
$ws.Columns.Item(8).ColumnWidth = 13.333333333333334
$ws.Columns.Item(12).ColumnWidth = 15.5

# ---------------------------------------------------------------------
# 4) Leave the selection where the author ended up (cell D9) instead of
#    the previous L7.
# ---------------------------------------------------------------------

$ws.Range("D9").Select() | Out-Null
